# Auto-generated Excel COM-interop edit script
# Applies the Typhon_Profits.xlsx leve-profit recalculation updates
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 313.9091
$ws.Range("I33").Value = 340.3
$ws.Range("J33").Value = 50
$ws.Range("K33").Value = 340.3
$ws.Range("L33").Value = 50
$ws.Range("M33").Value = -111.3
$ws.Range("N33").Value = -508
$ws.Range("H51").Value = 7500
$ws.Range("I51").Value = 7500
$ws.Range("K51").Value = 7500
$ws.Range("M51").Value = -7016
$ws.Range("H57").Value = 23162.5
$ws.Range("J57").Value = 23162.5
$ws.Range("L57").Value = 69487.5
$ws.Range("N57").Value = -70485.5
$ws.Range("H129").Value = 223355.22
$ws.Range("I129").Value = 283.33334
$ws.Range("J129").Value = 239288.92
$ws.Range("K129").Value = 850.0000200000001
$ws.Range("L129").Value = 717866.76
$ws.Range("M129").Value = 4149.99998
$ws.Range("N129").Value = -727866.76
$ws.Range("H132").Value = 3374.2068
$ws.Range("I132").Value = 3434.04
$ws.Range("J132").Value = 3000.25
$ws.Range("K132").Value = 10302.12
$ws.Range("L132").Value = 9000.75
$ws.Range("M132").Value = -7772.119999999999
$ws.Range("N132").Value = -14060.75
$ws.Range("H137").Value = 1989.5518
$ws.Range("I137").Value = 1833.2084
$ws.Range("J137").Value = 2740
$ws.Range("K137").Value = 5499.6252
$ws.Range("L137").Value = 8220
$ws.Range("M137").Value = -2949.6252
$ws.Range("N137").Value = -13320
$ws.Range("H138").Value = 10871805
$ws.Range("I138").Value = 22223090
$ws.Range("J138").Value = 3552.5532
$ws.Range("K138").Value = 66669270
$ws.Range("L138").Value = 10657.6596
$ws.Range("M138").Value = -66664130
$ws.Range("N138").Value = -20937.6596

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3926
$ws.Range("I32").Value = 2589.3538
$ws.Range("J32").Value = 17626.625
$ws.Range("K32").Value = 2589.3538
$ws.Range("L32").Value = 17626.625
$ws.Range("M32").Value = -2302.3538
$ws.Range("N32").Value = -18200.625
$ws.Range("H74").Value = 25642922
$ws.Range("I74").Value = 27779586
$ws.Range("J74").Value = 2966.3333
$ws.Range("K74").Value = 27779586
$ws.Range("L74").Value = 2966.3333
$ws.Range("M74").Value = -27778712
$ws.Range("N74").Value = -4714.3333
$ws.Range("H77").Value = 25642922
$ws.Range("I77").Value = 27779586
$ws.Range("J77").Value = 2966.3333
$ws.Range("K77").Value = 138897930
$ws.Range("L77").Value = 14831.6665
$ws.Range("M77").Value = -138893562
$ws.Range("N77").Value = -23567.6665
$ws.Range("H132").Value = 9346.229499999999
$ws.Range("I132").Value = 1096.7778
$ws.Range("J132").Value = 72984.86
$ws.Range("K132").Value = 3290.3334
$ws.Range("L132").Value = 218954.58
$ws.Range("M132").Value = -760.3334000000004
$ws.Range("N132").Value = -224014.58

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3364.976
$ws.Range("I134").Value = 3772.9688
$ws.Range("J134").Value = 2059.4
$ws.Range("K134").Value = 11318.9064
$ws.Range("L134").Value = 6178.200000000001
$ws.Range("M134").Value = -8783.9064
$ws.Range("N134").Value = -11248.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3137.2195
$ws.Range("I31").Value = 1767.129
$ws.Range("J31").Value = 7384.5
$ws.Range("K31").Value = 1767.129
$ws.Range("L31").Value = 7384.5
$ws.Range("M31").Value = -1472.129
$ws.Range("N31").Value = -7974.5
$ws.Range("H34").Value = 3137.2195
$ws.Range("I34").Value = 1767.129
$ws.Range("J34").Value = 7384.5
$ws.Range("K34").Value = 1767.129
$ws.Range("L34").Value = 7384.5
$ws.Range("M34").Value = -1565.129
$ws.Range("N34").Value = -7788.5
$ws.Range("H52").Value = 39994.5
$ws.Range("J52").Value = 39994.5
$ws.Range("L52").Value = 39994.5
$ws.Range("N52").Value = -40582.5
$ws.Range("H107").Value = 1446.2
$ws.Range("I107").Value = 1463.1666
$ws.Range("J107").Value = 1440.8422
$ws.Range("K107").Value = 1463.1666
$ws.Range("L107").Value = 1440.8422
$ws.Range("M107").Value = 456.8334
$ws.Range("N107").Value = -5280.8422
$ws.Range("H132").Value = 2409.625
$ws.Range("I132").Value = 1809.4839
$ws.Range("J132").Value = 21014
$ws.Range("K132").Value = 5428.4517
$ws.Range("L132").Value = 63042
$ws.Range("M132").Value = -2898.4517
$ws.Range("N132").Value = -68102
$ws.Range("H134").Value = 1193.6522
$ws.Range("I134").Value = 1023.9474
$ws.Range("J134").Value = 1999.75
$ws.Range("K134").Value = 3071.8422
$ws.Range("L134").Value = 5999.25
$ws.Range("M134").Value = -536.8422
$ws.Range("N134").Value = -11069.25
$ws.Range("H137").Value = 24990
$ws.Range("J137").Value = 24990
$ws.Range("L137").Value = 24990
$ws.Range("N137").Value = -35190

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 49.105263
$ws.Range("I12").Value = 7.857143
$ws.Range("J12").Value = 73.166664
$ws.Range("K12").Value = 23.571429
$ws.Range("L12").Value = 219.499992
$ws.Range("M12").Value = 149.428571
$ws.Range("N12").Value = -565.499992
$ws.Range("H62").Value = 6287.9165
$ws.Range("J62").Value = 7984.875
$ws.Range("L62").Value = 23954.625
$ws.Range("N62").Value = -25326.625
$ws.Range("H65").Value = 6287.9165
$ws.Range("J65").Value = 7984.875
$ws.Range("L65").Value = 71863.875
$ws.Range("N65").Value = -78727.875
$ws.Range("H68").Value = 900
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 900
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H131").Value = 758.48486
$ws.Range("J131").Value = 770.21277
$ws.Range("L131").Value = 2310.63831
$ws.Range("N131").Value = -12390.63831
$ws.Range("H132").Value = 633
$ws.Range("I132").Value = 450
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 4050
$ws.Range("L132").Value = 8991
$ws.Range("M132").Value = -1520
$ws.Range("N132").Value = -14051
$ws.Range("H134").Value = 2179.7
$ws.Range("J134").Value = 3600.625
$ws.Range("L134").Value = 10801.875
$ws.Range("N134").Value = -20941.875
$ws.Range("H136").Value = 1883.1765
$ws.Range("J136").Value = 4503.5
$ws.Range("L136").Value = 13510.5
$ws.Range("N136").Value = -23710.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 78433390
$ws.Range("I122").Value = 30304716
$ws.Range("J122").Value = 166669300
$ws.Range("K122").Value = 90914148
$ws.Range("L122").Value = 500007900
$ws.Range("M122").Value = -90911698
$ws.Range("N122").Value = -500012800
$ws.Range("H132").Value = 31615.666
$ws.Range("I132").Value = 4431.2
$ws.Range("J132").Value = 167538
$ws.Range("K132").Value = 13293.6
$ws.Range("L132").Value = 502614
$ws.Range("M132").Value = -10763.6
$ws.Range("N132").Value = -507674

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 43483244
$ws.Range("I7").Value = 100003310
$ws.Range("J7").Value = 6269.5386
$ws.Range("K7").Value = 100003310
$ws.Range("L7").Value = 6269.5386
$ws.Range("M7").Value = -100003198
$ws.Range("N7").Value = -6493.5386
$ws.Range("H93").Value = 1014.25
$ws.Range("I93").Value = 960.9167
$ws.Range("K93").Value = 960.9167
$ws.Range("M93").Value = 287.0833
$ws.Range("H126").Value = 43483244
$ws.Range("I126").Value = 100003310
$ws.Range("J126").Value = 6269.5386
$ws.Range("K126").Value = 300009930
$ws.Range("L126").Value = 18808.6158
$ws.Range("M126").Value = -300007460
$ws.Range("N126").Value = -23748.6158
$ws.Range("H132").Value = 1056.3529
$ws.Range("I132").Value = 1067.42
$ws.Range("K132").Value = 3202.26
$ws.Range("M132").Value = -672.2600000000002
$ws.Range("H136").Value = 828.902
$ws.Range("I136").Value = 767.8409
$ws.Range("J136").Value = 1212.7142
$ws.Range("K136").Value = 2303.5227
$ws.Range("L136").Value = 3638.1426
$ws.Range("M136").Value = 246.4773
$ws.Range("N136").Value = -8738.142599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 600
$ws.Range("I100").Value = 700
$ws.Range("J100").Value = 300
$ws.Range("K100").Value = 1400
$ws.Range("L100").Value = 600
$ws.Range("M100").Value = -859
$ws.Range("N100").Value = -1682
$ws.Range("H126").Value = 1079.7878
$ws.Range("I126").Value = 1079.7878
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3239.3634
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -769.3634000000002
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 700.55817
$ws.Range("I132").Value = 636.3611
$ws.Range("J132").Value = 1030.7142
$ws.Range("K132").Value = 1909.0833
$ws.Range("L132").Value = 3092.1426
$ws.Range("M132").Value = 620.9167000000002
$ws.Range("N132").Value = -8152.142599999999
$ws.Range("H136").Value = 18183952
$ws.Range("I136").Value = 25642026
$ws.Range("J136").Value = 4900.5625
$ws.Range("K136").Value = 4900.5625
$ws.Range("L136").Value = 14701.6875
$ws.Range("M136").Value = -76923528
$ws.Range("N136").Value = -19801.6875
